# Apply latest cryptocurrency market data to the "cryptos" worksheet.
# For rows 28/29, 32/33 and 41/42 the coin order also changed, so Coin (B)
# and Link (C) are updated together with Price (D) and Volume(1h) (E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, Coin(B), Link(C), Price(D), Volume1h(E)
# $null means "leave the existing value in that column untouched".
$updates = @(
    @{ Row = 2; B = $null; C = $null; D = '44.431.35'; E = '  +3.65%  ' },
    @{ Row = 3; B = $null; C = $null; D = '2.270.37'; E = '  +2.66%  ' },
    @{ Row = 4; B = $null; C = $null; D = $null; E = '  +0.03%  ' },
    @{ Row = 5; B = $null; C = $null; D = '320.45'; E = '  +1.36%  ' },
    @{ Row = 6; B = $null; C = $null; D = '105.24'; E = '  +5.76%  ' },
    @{ Row = 7; B = $null; C = $null; D = $null; E = '  +0.42%  ' },
    @{ Row = 8; B = $null; C = $null; D = $null; E = '  -0.13%  ' },
    @{ Row = 9; B = $null; C = $null; D = $null; E = '  +1.92%  ' },
    @{ Row = 10; B = $null; C = $null; D = '38.58'; E = '  +4.73%  ' },
    @{ Row = 11; B = $null; C = $null; D = '0.0844'; E = '  +1.90%  ' },
    @{ Row = 12; B = $null; C = $null; D = $null; E = '  +1.78%  ' },
    @{ Row = 13; B = $null; C = $null; D = $null; E = '  +0.67%  ' },
    @{ Row = 14; B = $null; C = $null; D = '0.882'; E = '  +2.70%  ' },
    @{ Row = 15; B = $null; C = $null; D = '2.616.57'; E = '  +2.82%  ' },
    @{ Row = 16; B = $null; C = $null; D = '14.54'; E = '  +2.48%  ' },
    @{ Row = 17; B = $null; C = $null; D = '2.274.82'; E = '  +2.91%  ' },
    @{ Row = 18; B = $null; C = $null; D = '44.295.65'; E = '  +3.63%  ' },
    @{ Row = 19; B = $null; C = $null; D = '13.97'; E = '  -4.78%  ' },
    @{ Row = 20; B = $null; C = $null; D = $null; E = '  +4.28%  ' },
    @{ Row = 21; B = $null; C = $null; D = $null; E = '  +1.96%  ' },
    @{ Row = 22; B = $null; C = $null; D = '66.30'; E = '  +1.56%  ' },
    @{ Row = 23; B = $null; C = $null; D = '3.19'; E = '  +1.35%  ' },
    @{ Row = 24; B = $null; C = $null; D = '239.41'; E = '  +1.37%  ' },
    @{ Row = 25; B = $null; C = $null; D = $null; E = '  +3.44%  ' },
    @{ Row = 26; B = $null; C = $null; D = $null; E = '  -0.24%  ' },
    @{ Row = 27; B = $null; C = $null; D = '10.19'; E = '  +1.73%  ' },
    @{ Row = 28; B = 'InjectiveProtocol'; C = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; D = '38.55'; E = '  +12.51%  ' },
    @{ Row = 29; B = 'Toncoin'; C = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D = '2.22'; E = '  +0.17%  ' },
    @{ Row = 30; B = $null; C = $null; D = $null; E = '  +2.38%  ' },
    @{ Row = 31; B = $null; C = $null; D = $null; E = '  +0.42%  ' },
    @{ Row = 32; B = 'Monero'; C = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D = '161.79'; E = '  +3.78%  ' },
    @{ Row = 33; B = 'Hedera'; C = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; D = '0.0884'; E = '  -1.25%  ' },
    @{ Row = 34; B = $null; C = $null; D = $null; E = '  -0.90%  ' },
    @{ Row = 35; B = $null; C = $null; D = $null; E = '  +9.29%  ' },
    @{ Row = 36; B = $null; C = $null; D = '2.02'; E = '  +4.07%  ' },
    @{ Row = 37; B = $null; C = $null; D = $null; E = '  -0.28%  ' },
    @{ Row = 38; B = $null; C = $null; D = '0.122'; E = '  +0.23%  ' },
    @{ Row = 39; B = $null; C = $null; D = $null; E = '  +1.86%  ' },
    @{ Row = 40; B = $null; C = $null; D = '4.43'; E = '  +0.23%  ' },
    @{ Row = 41; B = 'Celestia'; C = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'; D = '15.57'; E = '  +24.35%  ' },
    @{ Row = 42; B = 'VeChain'; C = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D = '0.0328'; E = '  +1.03%  ' },
    @{ Row = 43; B = $null; C = $null; D = $null; E = '  +0.31%  ' },
    @{ Row = 44; B = $null; C = $null; D = '1.769.41'; E = '  -6.55%  ' },
    @{ Row = 45; B = $null; C = $null; D = $null; E = '  +0.52%  ' },
    @{ Row = 46; B = $null; C = $null; D = '87.02'; E = '  -1.05%  ' },
    @{ Row = 47; B = $null; C = $null; D = '5.44'; E = '  +1.67%  ' },
    @{ Row = 48; B = $null; C = $null; D = '60.54'; E = '  -1.00%  ' },
    @{ Row = 49; B = $null; C = $null; D = '75.06'; E = '  -1.82%  ' },
    @{ Row = 50; B = $null; C = $null; D = '1.70'; E = '  +6.77%  ' },
    @{ Row = 51; B = $null; C = $null; D = '104.19'; E = '  +1.81%  ' }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.B) { $ws.Cells.Item($r, 2).Value = $u.B }   # Column B: Coin
    if ($null -ne $u.C) { $ws.Cells.Item($r, 3).Value = $u.C }   # Column C: Link
    if ($null -ne $u.D) {
        # Column D: Price - force text so values like "60.54" or "0.0844"
        # are not reinterpreted as numbers (which would drop trailing zeros
        # or collapse "44.431.35" style strings).
        $cellD = $ws.Cells.Item($r, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
        $cellD.ClearFormats()
    }
    if ($null -ne $u.E) { $ws.Cells.Item($r, 5).Value = $u.E }   # Column E: Volume(1h)
}

